# Daily auto-push update: a new hourly data point for 2026/01/10 (Sat)
# was recorded, inserted in its correct chronological position (row 618).
# Every later row shifts down by one; the last row (formerly 659) becomes
# row 660.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before row 618; rows 618-659 shift down to 619-660,
# and the sheet's dimension grows from D659 to D660 automatically.
$ws.Rows(618).Insert()

# Column A holds the date as literal text (e.g. "2026/01/10"), not a real
# date value, matching every other row in the sheet. Forcing the cell to
# Text format before assignment keeps Excel from auto-converting the
# yyyy/mm/dd-looking string into a date serial number; re-applying the
# neighboring cell's style afterwards drops the now-unneeded "@" number
# format so the cell's style matches the rest of the column (no explicit
# style override).
$ws.Range("A618").NumberFormat = "@"
$ws.Range("A618").Value = "2026/01/10"
$ws.Range("A618").Style = $ws.Range("A2").Style

$ws.Range("B618").Value = "土"
$ws.Range("C618").Value = 10
$ws.Range("D618").Value = 201
